# Closer, still need route param to work.
# Adds a new column G with an "error_log( print_R($<field>, TRUE ));" helper
# formula, flips the F "running concatenation" column to accumulate the
# other direction (F(n) = F(n-1) & E(n) instead of E(n) & F(n-1)), and
# re-tags a handful of Sheet2 VLOOKUP source rows from type "i" to a new
# type "x".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: re-tag the lookup rows that become the new "x" type ------
# (ConsentToPublicPictures, InstructorPaymentFree, include, InstructorFlag,
#  bdayinclude, sendWelcomeCard)
$ws2.Range("I20").Value = "x"
$ws2.Range("I22").Value = "x"
$ws2.Range("I24").Value = "x"
$ws2.Range("I25").Value = "x"
$ws2.Range("I30").Value = "x"
$ws2.Range("I41").Value = "x"

# --- Sheet1: flip column F to accumulate as F(n-1)&E(n) ---------------
$ws1.Range("F2").Formula = "=E2"
$ws1.Range("F3").Formula = "=F2&E3"
$ws1.Range("F4:F34").Formula = "=F3&E4"

# --- Sheet1: new column G with the error_log helper formula -----------
$ws1.Range("G2").Formula = '= "error_log( print_R($"&TRIM(A2)&", TRUE ));"'
$ws1.Range("G3:G33").Formula = '= "error_log( print_R($"&TRIM(A3)&", TRUE ));"'

# Column F is now wide enough to show the (now-short) running string, and
# the new column G needs a best-fit-ish width for its long formula result.
$ws1.Columns.Item(6).ColumnWidth = 29.33

# --- View state: scroll sheet1 so column B is left-most and select the
# new G2:G33 range (mirrors where the author was working) -------------
$ws1.Activate() | Out-Null
$ws1.Range("G2:G33").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2

# Sheet2's window was scrolled up a couple of rows too.
$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 14

# Leave Sheet1 as the active/selected sheet (matches tabSelected in the
# original file).
$ws1.Activate() | Out-Null
$ws1.Range("G2:G33").Select() | Out-Null
